$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F; this shifts the existing F:K
# columns (and their formatting) one column to the right, to G:L.
$ws.Range("F1").EntireColumn.Insert()

# Populate the newly inserted column F.
$ws.Range("F1").Value = "Unnamed: 0.1.1.1.1"
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 4

# Row 7 (the extra "crossectional randomforest" row) gets an index value
# in column B as well, matching the other data rows.
$ws.Range("B7").Value = 5
